# JLC_BOM.xlsx update
#
# The author shortened the multi-designator range labels in column B so the
# repeated letter prefix after the dash is dropped (e.g. "C7-C15" -> "C7-15",
# "D1-D64" -> "D1-64", "D65-D192" -> "D65-192", "D193-D199" -> "D193-199").
#
# Cell values are updated in the same order the author touched them so the
# shared-string table grows in the same sequence, then the active selection
# is left on D16 to match where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "C7-15"
$ws.Range("B9").Value  = "D65-192"
$ws.Range("B10").Value = "D193-199"
$ws.Range("B8").Value  = "D1-64"

$ws.Range("D16").Select()
